# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E23) is re-sorted from descending
# (2105,2104,2103,2102,2101,2012,2011,2010) to ascending
# (2010,2011,2012,2101,2102,2103,2104,2105), the "Salario Basico" column
# (G16:G23) is updated to the new value for every period, and the
# "Valor Mora" column (F16:F23) keeps its special value (28090) attached
# to period 2105, which is now the last row instead of the first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order with their Valor Mora (F) / Salario Basico (G)
$rows = @(
    @{ Row = 16; Periodo = "2010"; ValorMora = 35112; SalarioBasico = 877803 },
    @{ Row = 17; Periodo = "2011"; ValorMora = 35112; SalarioBasico = 877803 },
    @{ Row = 18; Periodo = "2012"; ValorMora = 35112; SalarioBasico = 877803 },
    @{ Row = 19; Periodo = "2101"; ValorMora = 35112; SalarioBasico = 877803 },
    @{ Row = 20; Periodo = "2102"; ValorMora = 35112; SalarioBasico = 877803 },
    @{ Row = 21; Periodo = "2103"; ValorMora = 35112; SalarioBasico = 877803 },
    @{ Row = 22; Periodo = "2104"; ValorMora = 35112; SalarioBasico = 877803 },
    @{ Row = 23; Periodo = "2105"; ValorMora = 28090; SalarioBasico = 877803 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("E$row").Value = $r.Periodo
    $ws.Range("F$row").Value = $r.ValorMora
    $ws.Range("G$row").Value = $r.SalarioBasico
}
